# Revert "Auto stash before revert of "Update to version 3.2.1""
# Sets the BAU Guaranteed Dispatch Percentage for "biomass" in year 2015
# back to 0 (the shared formulas across the row recompute automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")
$ws.Range("B9").Value = 0
